# Reorganise the folder and rename some files.
#
# The underlying XML diff is mostly Excel-version/authoring metadata
# (fileVersion, xr:* revision GUIDs, the author's local absPath, namespace
# bumps) that gets stamped automatically whenever the workbook is
# re-saved and isn't reachable through the Excel object model. The
# functionally meaningful, scriptable changes are:
#   - column A gets an explicit width; columns C/E/G/I get new widths
#   - the active selection moves from A21:XFD21 to J4
#   - the sheet view zoom changes to 84%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------
# Range.ColumnWidth is expressed in "characters" of the Normal style
# and Excel stores a slightly larger number in the OOXML <col width=.../>
# (characters + ~5/6). Back the COM value off by that fixed offset so the
# persisted XML width lands on the target value.
$widthOffset = 0.8333333333333334

$ws.Range("A1").EntireColumn.ColumnWidth = (9 - $widthOffset)
$ws.Range("C1").EntireColumn.ColumnWidth = (10.83203125 - $widthOffset)
$ws.Range("E1").EntireColumn.ColumnWidth = (10.83203125 - $widthOffset)
$ws.Range("G1").EntireColumn.ColumnWidth = (10.83203125 - $widthOffset)
$ws.Range("I1").EntireColumn.ColumnWidth = (10.5 - $widthOffset)

# --- Selection -------------------------------------------------------
$ws.Range("J4").Select() | Out-Null

# --- Zoom --------------------------------------------------------------
$excel.ActiveWindow.Zoom = 84
